$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test result status in C2 from "Passed" to "Failed"
$ws.Range("C2").Value = "Failed"

# Update the dates in H2 and H3 from "16/04/2021" to "17/04/2021"
$ws.Range("H2").Value = "17/04/2021"
$ws.Range("H3").Value = "17/04/2021"
